$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows 4-6, repeating the data from rows 2 and 3 (row2, row3, row2)
$ws.Range("A4").Value = "user4481660221383"
$ws.Range("B4").Value = " @K4XrpZotyHPt"

$ws.Range("A5").Value = "user25765839695551"
$ws.Range("B5").Value = " @K4I060Ckd395"

$ws.Range("A6").Value = "user4481660221383"
$ws.Range("B6").Value = " @K4XrpZotyHPt"

# Update the selected cell like the original edit session did
$ws.Range("B13").Select()

$wb.Save()
